$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A29").Value = "13 marras"
$ws.Range("C29").Value = "Köysidemon parantelua, "
$ws.Range("B29").Value = "10.00-11.30, 12.00-12.30"

$ws.Range("B29").NumberFormat = "h:mm"
$ws.Range("B29").WrapText = $true
$ws.Range("C29").WrapText = $true
$ws.Range("A29:G29").EntireRow.AutoFit()

$ws.Range("B30").Select()
